
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reg")

# --- Update the Username column (C) for each registration row ---
$ws.Range("C2").Value = "akhilbingi3213"
$ws.Range("C4").Value = "akhilbingi324455"
$ws.Range("C5").Value = "akhilbingi321566"
$ws.Range("C6").Value = "akhilbingi3212777"
$ws.Range("C7").Value = "nffmf88"
$ws.Range("C8").Value = "akhilhdiw599"
$ws.Range("C9").Value = "akhilhdingi900"

# C3 additionally gets a mailto hyperlink whose visible text differs from the
# cell's stored value, so add the hyperlink (which sets both) first, then
# overwrite the cell text with the real username afterwards.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:akhilbingi6423@gmail.com", "", "", "akhilbingi6423@gmail.com")
$ws.Range("C3").Value = "akhilbingi6423444"

# --- Widen column C (username) and drop its autofit flag ---
$ws.Columns.Item(3).ColumnWidth = 36.83

# --- Update the active selection on the sheet ---
$ws.Activate() | Out-Null
$ws.Range("C9").Select() | Out-Null
